$d = $word.ActiveDocument

$pairs = @(
    @("2024-05-23 Thursday", "2024-05-24 Friday"),
    @("976÷7=139, 3", "113÷9=12, 5"),
    @("510÷4=127, 2", "344÷6=57, 2"),
    @("582÷5=116, 2", "724÷7=103, 3"),
    @("850÷6=141, 4", "960÷3=320, 0"),
    @("596÷8=74, 4", "204÷8=25, 4"),
    @("542÷4=135, 2", "823÷2=411, 1"),
    @("651÷5=130, 1", "216÷2=108, 0"),
    @("397÷2=198, 1", "342÷8=42, 6"),
    @("851÷9=94, 5", "809÷5=161, 4"),
    @("483÷5=96, 3", "639÷7=91, 2"),
    @("408÷4=102, 0", "556÷4=139, 0"),
    @("506÷3=168, 2", "903÷2=451, 1"),
    @("600÷2=300, 0", "394÷8=49, 2"),
    @("110÷5=22, 0", "410÷3=136, 2"),
    @("913÷6=152, 1", "506÷8=63, 2"),
    @("898÷9=99, 7", "163÷4=40, 3"),
    @("595÷3=198, 1", "452÷4=113, 0"),
    @("380÷3=126, 2", "745÷5=149, 0"),
    @("137÷2=68, 1", "531÷2=265, 1"),
    @("612÷6=102, 0", "484÷3=161, 1"),
    @("140÷3=46, 2", "563÷2=281, 1"),
    @("424÷8=53, 0", "267÷8=33, 3"),
    @("178÷5=35, 3", "214÷6=35, 4"),
    @("176÷9=19, 5", "168÷9=18, 6"),
    @("262÷5=52, 2", "748÷2=374, 0")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
